# Update "paises.xlsx" COVID country stats (Pais sheet).
# The source data feed refreshed a handful of country rows and the sheet
# stays sorted by "Casos totales" (col B) descending, so a few entries
# change rank and swap places with their neighbours:
#   - Dinamarca overtakes Noruega (rows 32/33)
#   - Banglades surges past 8 countries (row 73 -> row 65, others shift down)
#   - Consejo Danes para los Refugiados overtakes Isla de Man (rows 115/116)
# Each affected row is written directly with its final country name + stats,
# and the "Datos actualizados" timestamp cell is bumped to 10:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp (A1): 09:52 -> 10:52
$ws.Range("A1").Value = "Datos actualizados a 16 de Abril de 2020 a las 10:52"

# Row 18: Suiza
$ws.Cells.Item(18, 1).Value = "Suiza"
$ws.Cells.Item(18, 2).Value = 26336
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 15400
$ws.Cells.Item(18, 5).Value = 9687
$ws.Cells.Item(18, 6).Value = 386
$ws.Cells.Item(18, 7).Value = 10
$ws.Cells.Item(18, 8).Value = 1249

# Row 20: Austria
$ws.Cells.Item(20, 1).Value = "Austria"
$ws.Cells.Item(20, 2).Value = 14404
$ws.Cells.Item(20, 3).Value = 54
$ws.Cells.Item(20, 4).Value = 8098
$ws.Cells.Item(20, 5).Value = 5913
$ws.Cells.Item(20, 6).Value = 232
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 393

# Row 30: Polonia
$ws.Cells.Item(30, 1).Value = "Polonia"
$ws.Cells.Item(30, 2).Value = 7771
$ws.Cells.Item(30, 3).Value = 189
$ws.Cells.Item(30, 4).Value = 774
$ws.Cells.Item(30, 5).Value = 6705
$ws.Cells.Item(30, 6).Value = 160
$ws.Cells.Item(30, 7).Value = 6
$ws.Cells.Item(30, 8).Value = 292

# Row 31: Rumania
$ws.Cells.Item(31, 1).Value = "Rumania"
$ws.Cells.Item(31, 2).Value = 7216
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 1217
$ws.Cells.Item(31, 5).Value = 5612
$ws.Cells.Item(31, 6).Value = 245
$ws.Cells.Item(31, 7).Value = 15
$ws.Cells.Item(31, 8).Value = 387

# Row 32: Dinamarca
$ws.Cells.Item(32, 1).Value = "Dinamarca"
$ws.Cells.Item(32, 2).Value = 6879
$ws.Cells.Item(32, 3).Value = 198
$ws.Cells.Item(32, 4).Value = 2748
$ws.Cells.Item(32, 5).Value = 3822
$ws.Cells.Item(32, 6).Value = 89
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 309

# Row 33: Noruega
$ws.Cells.Item(33, 1).Value = "Noruega"
$ws.Cells.Item(33, 2).Value = 6798
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = 32
$ws.Cells.Item(33, 5).Value = 6616
$ws.Cells.Item(33, 6).Value = 64
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 150

# Row 35: Australia
$ws.Cells.Item(35, 1).Value = "Australia"
$ws.Cells.Item(35, 2).Value = 6468
$ws.Cells.Item(35, 3).Value = 21
$ws.Cells.Item(35, 4).Value = 3747
$ws.Cells.Item(35, 5).Value = 2658
$ws.Cells.Item(35, 6).Value = 66
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 63

# Row 39: Filipinas
$ws.Cells.Item(39, 1).Value = "Filipinas"
$ws.Cells.Item(39, 2).Value = 5660
$ws.Cells.Item(39, 3).Value = 207
$ws.Cells.Item(39, 4).Value = 435
$ws.Cells.Item(39, 5).Value = 4863
$ws.Cells.Item(39, 6).Value = 1
$ws.Cells.Item(39, 7).Value = 13
$ws.Cells.Item(39, 8).Value = 362

# Row 59: Moldavia
$ws.Cells.Item(59, 1).Value = "Moldavia"
$ws.Cells.Item(59, 2).Value = 2049
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 171
$ws.Cells.Item(59, 5).Value = 1827
$ws.Cells.Item(59, 6).Value = 80
$ws.Cells.Item(59, 7).Value = 5
$ws.Cells.Item(59, 8).Value = 51

# Row 63: Barein
$ws.Cells.Item(63, 1).Value = "Barein"
$ws.Cells.Item(63, 2).Value = 1673
$ws.Cells.Item(63, 3).Value = 2
$ws.Cells.Item(63, 4).Value = 663
$ws.Cells.Item(63, 5).Value = 1003
$ws.Cells.Item(63, 6).Value = 3
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 7

# Row 65: Banglades
$ws.Cells.Item(65, 1).Value = "Banglades"
$ws.Cells.Item(65, 2).Value = 1572
$ws.Cells.Item(65, 3).Value = 341
$ws.Cells.Item(65, 4).Value = 49
$ws.Cells.Item(65, 5).Value = 1463
$ws.Cells.Item(65, 6).Value = 1
$ws.Cells.Item(65, 7).Value = 10
$ws.Cells.Item(65, 8).Value = 60

# Row 66: Irak
$ws.Cells.Item(66, 1).Value = "Irak"
$ws.Cells.Item(66, 2).Value = 1415
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 4).Value = 812
$ws.Cells.Item(66, 5).Value = 524
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 79

# Row 67: Kuwait
$ws.Cells.Item(67, 1).Value = "Kuwait"
$ws.Cells.Item(67, 2).Value = 1405
$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(67, 4).Value = 206
$ws.Cells.Item(67, 5).Value = 1196
$ws.Cells.Item(67, 6).Value = 31
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 3

# Row 68: Nueva Zelanda
$ws.Cells.Item(68, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(68, 2).Value = 1401
$ws.Cells.Item(68, 3).Value = 15
$ws.Cells.Item(68, 4).Value = 770
$ws.Cells.Item(68, 5).Value = 622
$ws.Cells.Item(68, 6).Value = 3
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 9

# Row 69: Estonia
$ws.Cells.Item(69, 1).Value = "Estonia"
$ws.Cells.Item(69, 2).Value = 1400
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 4).Value = 117
$ws.Cells.Item(69, 5).Value = 1248
$ws.Cells.Item(69, 6).Value = 10
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 35

# Row 70: Uzbekistan
$ws.Cells.Item(70, 1).Value = "Uzbekistan"
$ws.Cells.Item(70, 2).Value = 1349
$ws.Cells.Item(70, 3).Value = 47
$ws.Cells.Item(70, 4).Value = 107
$ws.Cells.Item(70, 5).Value = 1238
$ws.Cells.Item(70, 6).Value = 8
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 4

# Row 71: Kazajistan
$ws.Cells.Item(71, 1).Value = "Kazajistan"
$ws.Cells.Item(71, 2).Value = 1341
$ws.Cells.Item(71, 3).Value = 46
$ws.Cells.Item(71, 4).Value = 263
$ws.Cells.Item(71, 5).Value = 1062
$ws.Cells.Item(71, 6).Value = 22
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 16

# Row 72: Azerbaiyan
$ws.Cells.Item(72, 1).Value = "Azerbaiyan"
$ws.Cells.Item(72, 2).Value = 1253
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(72, 4).Value = 404
$ws.Cells.Item(72, 5).Value = 836
$ws.Cells.Item(72, 6).Value = 24
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 13

# Row 73: Eslovenia
$ws.Cells.Item(73, 1).Value = "Eslovenia"
$ws.Cells.Item(73, 2).Value = 1248
$ws.Cells.Item(73, 3).Value = 0
$ws.Cells.Item(73, 4).Value = 165
$ws.Cells.Item(73, 5).Value = 1022
$ws.Cells.Item(73, 6).Value = 34
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 61

# Row 78: Hong Kong
$ws.Cells.Item(78, 1).Value = "Hong Kong"
$ws.Cells.Item(78, 2).Value = 1018
$ws.Cells.Item(78, 3).Value = 1
$ws.Cells.Item(78, 4).Value = 485
$ws.Cells.Item(78, 5).Value = 529
$ws.Cells.Item(78, 6).Value = 9
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 4

# Row 115: Consejo Danes para los Refugiados
$ws.Cells.Item(115, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(115, 2).Value = 267
$ws.Cells.Item(115, 3).Value = 13
$ws.Cells.Item(115, 4).Value = 23
$ws.Cells.Item(115, 5).Value = 222
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = 22

# Row 116: Isla de Man
$ws.Cells.Item(116, 1).Value = "Isla de Man"
$ws.Cells.Item(116, 2).Value = 258
$ws.Cells.Item(116, 3).Value = 2
$ws.Cells.Item(116, 4).Value = 151
$ws.Cells.Item(116, 5).Value = 103
$ws.Cells.Item(116, 6).Value = 13
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 4

# Row 171: Republica del Chad
$ws.Cells.Item(171, 1).Value = "Republica del Chad"
$ws.Cells.Item(171, 2).Value = 27
$ws.Cells.Item(171, 3).Value = 4
$ws.Cells.Item(171, 4).Value = 5
$ws.Cells.Item(171, 5).Value = 22
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 0

# Row 175: Laos
$ws.Cells.Item(175, 1).Value = "Laos"
$ws.Cells.Item(175, 2).Value = 19
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 2
$ws.Cells.Item(175, 5).Value = 17
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0
